$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:2").Delete()
$ws.Range("A11").Value = "392c340c-6104-44fc-8759-e9565bc64ed9"
$ws.Range("B11").Value = "Burger"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "6.7"
$ws.Range("D11").Value = "NTU"
$ws.Range("E11").Value = "Burger"
$ws.Range("F11").Value = "Hot burger with sauce"
$ws.Range("A11:F11").ClearFormats()
